## This script applies the edit described by the commit:
## "[ADDITIONAL SCRAPING] added code to scrape more data about a player's
## batting performance in a match, also updated the excel sheets"
##
## Concretely:
##   1. Insert a brand-new "Player Info" worksheet as the first sheet,
##      with an ID/NAME/BATTING_HAND/BOWL_STYLE table for player 4104.
##   2. On "ODI Batting": rename header MATCH_CARD_LINK -> MATCH_CODE and
##      replace the full howstat.com scorecard URLs with the bare
##      MatchCode values.
##   3. On "ODI Bowling": same MATCH_CARD_LINK -> MATCH_CODE rename and
##      URL -> bare-code replacement.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update "ODI Batting" (MATCH_CARD_LINK column D)
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")

$batting.Range("D1").Value = "MATCH_CODE"

$batting.Range("D2:D6").NumberFormat = "@"
$batting.Range("D2").Value = "4402"
$batting.Range("D3").Value = "4406"
$batting.Range("D4").Value = "4410"
$batting.Range("D5").Value = "4435"
$batting.Range("D6").Value = "4436"

# ---------------------------------------------------------------------
# 2. Update "ODI Bowling" (MATCH_CARD_LINK column B)
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")

$bowling.Range("B1").Value = "MATCH_CODE"

$bowling.Range("B2").NumberFormat = "@"
$bowling.Range("B2").Value = "4436"

# ---------------------------------------------------------------------
# 3. Insert the new "Player Info" sheet in front of "ODI Batting"
# ---------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add($batting)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$hdr = $playerInfo.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4104"
$playerInfo.Range("B2").Value = "Mayank Anurag Agarwal"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Does Not Bowl | Unknown"

Write-Output "Player Info sheet added; MATCH_CARD_LINK columns converted to MATCH_CODE on ODI Batting/Bowling."
